$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on columns D and E data cells so values like
# '25.80' or percentages keep their exact textual representation instead
# of being auto-converted to numbers/losing trailing zeros.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '47.757.82'
$ws.Range('E2').Value = '  +0.37%  '
$ws.Range('D3').Value = '2.495.75'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '322.04'
$ws.Range('E5').Value = '  -0.20%  '
$ws.Range('D6').Value = '109.18'
$ws.Range('E6').Value = '  +3.85%  '
$ws.Range('E7').Value = '  -0.74%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '0.543'
$ws.Range('E9').Value = '  +0.11%  '
$ws.Range('D10').Value = '39.46'
$ws.Range('E10').Value = '  +3.46%  '
$ws.Range('D11').Value = '0.0811'
$ws.Range('E11').Value = '  -0.40%  '
$ws.Range('E12').Value = '  +0.53%  '
$ws.Range('D13').Value = '18.63'
$ws.Range('E13').Value = '  +1.69%  '
$ws.Range('D14').Value = '7.22'
$ws.Range('E14').Value = '  +0.72%  '
$ws.Range('D15').Value = '2.884.51'
$ws.Range('E15').Value = '  +0.11%  '
$ws.Range('D16').Value = '2.495.38'
$ws.Range('E16').Value = '  +0.24%  '
$ws.Range('D17').Value = '0.848'
$ws.Range('E17').Value = '  -0.14%  '
$ws.Range('D18').Value = '47.547.23'
$ws.Range('E18').Value = '  +0.18%  '
$ws.Range('D19').Value = '13.32'
$ws.Range('E19').Value = '  +4.04%  '
$ws.Range('D20').Value = '6.65'
$ws.Range('E20').Value = '  +0.83%  '
$ws.Range('D21').Value = '0.0₃0943'
$ws.Range('E21').Value = '  +0.55%  '
$ws.Range('D22').Value = '2.75'
$ws.Range('E22').Value = '  +14.43%  '
$ws.Range('D23').Value = '70.67'
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').Value = '247.44'
$ws.Range('E24').Value = '  -1.59%  '
$ws.Range('D25').Value = '2.56'
$ws.Range('E25').Value = '  -0.37%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').Value = '25.80'
$ws.Range('E27').Value = '  -1.48%  '
$ws.Range('D28').Value = '2.24'
$ws.Range('E28').Value = '  +1.18%  '
$ws.Range('E29').Value = '  -0.27%  '
$ws.Range('D30').Value = '0.138'
$ws.Range('E30').Value = '  +2.77%  '
$ws.Range('D31').Value = '34.79'
$ws.Range('E31').Value = '  -1.05%  '
$ws.Range('D32').Value = '49.87'
$ws.Range('E32').Value = '  +0.94%  '
$ws.Range('D33').Value = '20.25'
$ws.Range('E33').Value = '  +2.27%  '
$ws.Range('D34').Value = '5.34'
$ws.Range('E34').Value = '  -0.41%  '
$ws.Range('D35').Value = '0.0788'
$ws.Range('E35').Value = '  +0.46%  '
$ws.Range('E36').Value = '  +0.11%  '
$ws.Range('E37').Value = '  +2.02%  '
$ws.Range('E38').Value = '  -0.75%  '
$ws.Range('E39').Value = '  -1.51%  '
$ws.Range('E40').Value = '  +0.23%  '
$ws.Range('B41').Value = 'WEMIXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D41').Value = '2.22'
$ws.Range('E41').Value = '  -1.77%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').Value = '22.16'
$ws.Range('E42').Value = '  +3.84%  '
$ws.Range('D43').Value = '119.78'
$ws.Range('E43').Value = '  -1.95%  '
$ws.Range('E44').Value = '  +0.11%  '
$ws.Range('D45').Value = '1.994.72'
$ws.Range('E45').Value = '  +1.38%  '
$ws.Range('D46').Value = '3.04'
$ws.Range('E46').Value = '  +1.84%  '
$ws.Range('E48').Value = '  -0.40%  '
$ws.Range('E49').Value = '  -1.59%  '
$ws.Range('D50').Value = '5.23'
$ws.Range('E50').Value = '  -0.70%  '
$ws.Range('D51').Value = '56.68'
$ws.Range('E51').Value = '  +3.39%  '
